$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell text values first, in top-to-bottom row order, so that
# the resulting shared-strings table is built in the same order as a user
# typing the rows from top to bottom.
$ws.Range("A6").Value = "login"
$ws.Range("F6").Value = '//a[@class="_1_3w1N"]'
$ws.Range("A7").Value = "more"
$ws.Range("F7").Value = '//*[@id="container"]/div/div[2]/div[1]/div[2]/div[5]/div/div/div/div'
$ws.Range("A9").Value = "explore Plus"
$ws.Range("F9").Value = '//a[text()="Explore"]'
$ws.Range("A10").Value = "grocery"
$ws.Range("F10").Value = '//div[@class="xtXmba"and text()="Grocery"]'
$ws.Range("A12").Value = "view all for best electronics"
$ws.Range("F12").Value = "/html/body/div/div/div[4]/div[3]/div[1]/div/div[1]/div/div/a-absolute x path"
$ws.Range("A14").Value = "view all for best electronics"
$ws.Range("F14").Value = "//div[3]/div[1]/div/div[1]/div/div/a-relative x path"

# --- Rows 12-13: "view all for best electronics" -> absolute x path ---
$ws.Range("A12:E13").Merge()
$ws.Range("F12:I13").Merge()
$ws.Range("A12:E13").HorizontalAlignment = -4108
$ws.Range("F12:I13").HorizontalAlignment = -4108
$ws.Range("F12:I13").WrapText = $true

# --- Rows 14-15: "view all for best electronics" -> relative x path ---
$ws.Range("A14:E15").Merge()
$ws.Range("F14:I15").Merge()
$ws.Range("A14:E15").HorizontalAlignment = -4108
$ws.Range("F14:I15").HorizontalAlignment = -4108
$ws.Range("F14:I15").WrapText = $true

# --- Row 9: "explore Plus" -> Explore ---
$ws.Range("A9:E9").Merge()
$ws.Range("F9:I9").Merge()
$ws.Range("A9:E9").HorizontalAlignment = -4108
$ws.Range("F9:I9").HorizontalAlignment = -4108

# --- Rows 10-11: "grocery" -> Grocery ---
$ws.Range("A10:E11").Merge()
$ws.Range("F10:I11").Merge()
$ws.Range("A10:E11").HorizontalAlignment = -4108
$ws.Range("F10:I11").HorizontalAlignment = -4108
$ws.Range("F10:I11").WrapText = $true

# --- Row 6: "login" ---
$ws.Range("A6:E6").Merge()
$ws.Range("F6:I6").Merge()
$ws.Range("A6:E6").HorizontalAlignment = -4108
$ws.Range("F6:I6").HorizontalAlignment = -4108

# --- Rows 7-8: "more" ---
$ws.Range("A7:E8").Merge()
$ws.Range("F7:I8").Merge()
$ws.Range("A7:E8").HorizontalAlignment = -4108
$ws.Range("F7:I8").HorizontalAlignment = -4108
$ws.Range("F7:I8").WrapText = $true

$excel.CutCopyMode = $false
$ws.Range("F14:I15").Select()
